$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates that apply identically to both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same data).
$updates = @{
    2  = 8025
    3  = 7651
    5  = 188
    13 = 116
    14 = 1229
    16 = 45
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
